# Update the workbook to reflect a re-run of the SPADE script (new Start/End
# times) and updated R package version numbers, and drop the "backports"
# package row from the "Loaded only" packages table on the sessionInfo sheet.

$wb = $excel.ActiveWorkbook

# --- Info sheet: refresh the run's Start_time / End_time ---------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("B26").Value = "Thu Nov 19 15:23:38 2020"
$wsInfo.Range("B27").Value = "Thu Nov 19 15:23:45 2020"

# --- sessionInfo sheet: package version bumps + removed package --------
$wsSession = $wb.Worksheets.Item("sessionInfo")

# here: 0.1 -> 1.0.0
$wsSession.Range("G2").Value = "1.0.0"

# magrittr: 1.5 -> 2.0.1
$wsSession.Range("J3").Value = "2.0.1"

# rprojroot: 1.3-2 -> 2.0.2
$wsSession.Range("J10").Value = "2.0.2"

# Remove the "backports" / "1.1.10" row from the Loaded-only packages
# table (columns I/J), shifting the following "boot" / "1.3-25" row up
# and leaving the last row of that column pair empty.
$bootName = $wsSession.Range("I16").Value2
$bootVersion = $wsSession.Range("J16").Value2
$wsSession.Range("I15").Value = $bootName
$wsSession.Range("J15").Value = $bootVersion
$wsSession.Range("I16").ClearContents()
$wsSession.Range("J16").ClearContents()
